$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Spencer Wayne"
$ws.Range("I2").Value = "Spencer"
$ws.Range("K2").Value = "Wayne"

$ws.Range("A3").Select()
